$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 230.25
$ws.Range("I8").Value = 230.25
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 690.75
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -551.75
$ws.Range("N8").ClearContents()

$ws.Range("H15").Value = 1778.5883
$ws.Range("I15").Value = 1778.5883
$ws.Range("K15").Value = 5335.7649
$ws.Range("M15").Value = -5166.7649

$ws.Range("H61").Value = 351
$ws.Range("I61").Value = 151.5
$ws.Range("J61").Value = 750
$ws.Range("K61").Value = 454.5
$ws.Range("L61").Value = 2250
$ws.Range("M61").Value = -282.5
$ws.Range("N61").Value = -2594

$ws.Range("H99").Value = 2368.75
$ws.Range("I99").Value = 237.5
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 712.5
$ws.Range("L99").Value = 13500
$ws.Range("M99").Value = 785.5
$ws.Range("N99").Value = -16496

$ws.Range("H116").Value = 2897.5454
$ws.Range("I116").Value = 2957.1428
$ws.Range("J116").Value = 2793.25
$ws.Range("K116").Value = 2957.1428
$ws.Range("L116").Value = 2793.25
$ws.Range("M116").Value = 484.8571999999999
$ws.Range("N116").Value = -9677.25

$ws.Range("H132").Value = 1035.6786
$ws.Range("I132").Value = 961.5
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2884.5
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -354.5
$ws.Range("N132").Value = -11060

$ws.Range("H134").Value = 67725
$ws.Range("J134").Value = 67725
$ws.Range("L134").Value = 67725
$ws.Range("N134").Value = -77865

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9990.375
$ws.Range("I61").Value = 6573.5454
$ws.Range("J61").Value = 17507.4
$ws.Range("K61").Value = 6573.5454
$ws.Range("L61").Value = 17507.4
$ws.Range("M61").Value = -6361.5454
$ws.Range("N61").Value = -17931.4

$ws.Range("H132").Value = 2369.8838
$ws.Range("I132").Value = 1958.7778
$ws.Range("J132").Value = 4484.143
$ws.Range("K132").Value = 5876.3334
$ws.Range("L132").Value = 13452.429
$ws.Range("M132").Value = -3346.3334
$ws.Range("N132").Value = -18512.429

$ws.Range("H136").Value = 9990.375
$ws.Range("I136").Value = 6573.5454
$ws.Range("J136").Value = 17507.4
$ws.Range("K136").Value = 19720.6362
$ws.Range("L136").Value = 52522.2
$ws.Range("M136").Value = -17170.6362
$ws.Range("N136").Value = -57622.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 59000
$ws.Range("J59").Value = 59000
$ws.Range("L59").Value = 59000
$ws.Range("N59").Value = -60694

$ws.Range("H94").Value = 2094.12
$ws.Range("I94").Value = 2103.842
$ws.Range("J94").Value = 2063.3333
$ws.Range("K94").Value = 2103.842
$ws.Range("L94").Value = 2063.3333
$ws.Range("M94").Value = -1652.842
$ws.Range("N94").Value = -2965.3333

$ws.Range("H105").Value = 3984.4
$ws.Range("I105").Value = 4202.8823
$ws.Range("J105").Value = 3520.125
$ws.Range("K105").Value = 4202.8823
$ws.Range("L105").Value = 3520.125
$ws.Range("M105").Value = -2455.8823
$ws.Range("N105").Value = -7014.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1006
$ws.Range("I94").Value = 1012
$ws.Range("K94").Value = 1012
$ws.Range("M94").Value = -561

$ws.Range("H132").Value = 3082.92
$ws.Range("I132").Value = 2947.4358
$ws.Range("J132").Value = 3563.2727
$ws.Range("K132").Value = 8842.307400000002
$ws.Range("L132").Value = 10689.8181
$ws.Range("M132").Value = -6312.307400000002
$ws.Range("N132").Value = -15749.8181

$ws.Range("H134").Value = 2412.8286
$ws.Range("I134").Value = 1939.2693
$ws.Range("J134").Value = 3780.889
$ws.Range("K134").Value = 5817.8079
$ws.Range("L134").Value = 11342.667
$ws.Range("M134").Value = -3282.8079
$ws.Range("N134").Value = -16412.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 38774.457
$ws.Range("I14").Value = 38774.457
$ws.Range("K14").Value = 116323.371
$ws.Range("M14").Value = -116150.371

$ws.Range("H16").Value = 1225
$ws.Range("I16").Value = 633.3333
$ws.Range("K16").Value = 1899.9999
$ws.Range("M16").Value = -1726.9999

$ws.Range("H122").Value = 814.2857
$ws.Range("I122").Value = 225.125
$ws.Range("J122").Value = 1599.8334
$ws.Range("K122").Value = 2026.125
$ws.Range("L122").Value = 14398.5006
$ws.Range("M122").Value = 423.875
$ws.Range("N122").Value = -19298.5006

$ws.Range("H131").Value = 18193.678
$ws.Range("J131").Value = 21942.521
$ws.Range("L131").Value = 65827.56299999999
$ws.Range("N131").Value = -75907.56299999999

$ws.Range("H137").Value = 41551.69
$ws.Range("I137").Value = 1059.6666
$ws.Range("J137").Value = 53699.3
$ws.Range("K137").Value = 3178.9998
$ws.Range("L137").Value = 161097.9
$ws.Range("M137").Value = 1921.0002
$ws.Range("N137").Value = -171297.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 5007503
$ws.Range("I20").Value = 10000000
$ws.Range("K20").Value = 10000000
$ws.Range("M20").Value = -9999755

$ws.Range("H103").Value = 167500
$ws.Range("J103").Value = 167500
$ws.Range("L103").Value = 167500
$ws.Range("N103").Value = -169844

$ws.Range("H126").Value = 2711.4
$ws.Range("I126").Value = 1971.3846
$ws.Range("J126").Value = 4085.7144
$ws.Range("K126").Value = 5914.1538
$ws.Range("L126").Value = 12257.1432
$ws.Range("M126").Value = -3444.1538
$ws.Range("N126").Value = -17197.1432

$ws.Range("H132").Value = 6092.222
$ws.Range("I132").Value = 1890.625
$ws.Range("J132").Value = 12203.637
$ws.Range("K132").Value = 5671.875
$ws.Range("L132").Value = 36610.911
$ws.Range("M132").Value = -3141.875
$ws.Range("N132").Value = -41670.911

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 30000
$ws.Range("J24").Value = 30000
$ws.Range("L24").Value = 30000
$ws.Range("N24").Value = -30686

$ws.Range("H61").Value = 30998.545
$ws.Range("I61").Value = 30998.545
$ws.Range("K61").Value = 30998.545
$ws.Range("M61").Value = -30796.545

$ws.Range("H113").Value = 30998.545
$ws.Range("I113").Value = 30998.545
$ws.Range("K113").Value = 30998.545
$ws.Range("M113").Value = -28828.545

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 26583.334
$ws.Range("J20").Value = 26583.334
$ws.Range("L20").Value = 26583.334
$ws.Range("N20").Value = -27063.334

$ws.Range("H132").Value = 2741.9285
$ws.Range("I132").Value = 2857.682
$ws.Range("J132").Value = 2317.5
$ws.Range("K132").Value = 8573.045999999998
$ws.Range("L132").Value = 6952.5
$ws.Range("M132").Value = -6043.045999999998
$ws.Range("N132").Value = -12012.5
